# Sync updated NextBus feed values into the "NextBus1" sheet.
# (Mirrors an automated "Sync file from Google Drive" refresh: the
# EstimatedTimeOfArrival timestamps and MinutesToArrival counters move
# forward, one bus's TypeOfBus got corrected, and everything else on
# the sheet is left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NextBus1")

# Row 2 - Bishan Int / bus 52
$ws.Range("F2").Value = 45688.60689814815
$ws.Range("O2").Value = 13

# Row 3 - Gali Batu Ter / bus 184
$ws.Range("F3").Value = 45688.6058449074
$ws.Range("O3").Value = 11

# Row 4 - Gali Batu Ter / bus 75
$ws.Range("O4").Value = 23

# Row 5 - Eunos Int / bus 154
$ws.Range("F5").Value = 45688.60388888889
$ws.Range("O5").Value = 8

# Row 6 - Bt Batok Int / bus 61
$ws.Range("F6").Value = 45688.62048611111
$ws.Range("O6").Value = 32

# Row 7 - Hougang Ctrl Int / bus 151
$ws.Range("F7").Value = 45688.62586805555
$ws.Range("L7").Value = "DD"
$ws.Range("O7").Value = 40

# Row 8 - Hougang Ctrl Int / bus 74
$ws.Range("F8").Value = 45688.60770833334
$ws.Range("O8").Value = 14
